$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.129.42'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '3.541.65'
$ws.Range("E3").Value = '  +4.71%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''598.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.76%  '
$ws.Range("D6").Value = '''138.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.83%  '
$ws.Range("D7").Value = '3.542.69'
$ws.Range("E7").Value = '  +4.72%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("E10").Value = '  +4.33%  '
$ws.Range("D11").Value = '''6.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("E12").Value = '  +4.75%  '
$ws.Range("D13").Value = '4.144.96'
$ws.Range("E13").Value = '  +4.43%  '
$ws.Range("E14").Value = '  +4.64%  '
$ws.Range("D15").Value = '''27.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.14%  '
$ws.Range("D16").Value = '3.542.53'
$ws.Range("E16").Value = '  +3.59%  '
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").Value = '65.113.55'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = '''10.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.99%  '
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("E21").Value = '  +6.13%  '
$ws.Range("D22").Value = '''392.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.87%  '
$ws.Range("E23").Value = '  +5.19%  '
$ws.Range("D24").Value = '3.684.31'
$ws.Range("E24").Value = '  +4.40%  '
$ws.Range("D25").Value = '''73.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.26%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '''0.0000113'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.33%  '
$ws.Range("D28").Value = '''7.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.80%  '
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +4.85%  '
$ws.Range("D31").Value = '''8.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.68%  '
$ws.Range("D32").Value = '3.562.50'
$ws.Range("E32").Value = '  +4.66%  '
$ws.Range("D33").Value = '''1.39'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +23.96%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("E35").Value = '  +5.39%  '
$ws.Range("D36").Value = '''0.146'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.66%  '
$ws.Range("E37").Value = '  +10.66%  '
$ws.Range("D38").Value = '''170.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("E39").Value = '  +5.67%  '
$ws.Range("D40").Value = '''5.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.12%  '
$ws.Range("E41").Value = '  +8.62%  '
$ws.Range("E42").Value = '  +2.37%  '
$ws.Range("D43").Value = '''26.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +22.70%  '
$ws.Range("D44").Value = '''42.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("E46").Value = '  +4.06%  '
$ws.Range("E47").Value = '  +11.54%  '
$ws.Range("E48").Value = '  +6.31%  '
$ws.Range("E49").Value = '  +7.08%  '
$ws.Range("D50").Value = '2.403.68'
$ws.Range("E50").Value = '  +12.35%  '
$ws.Range("D51").Value = '''311.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +19.47%  '
